{"js": "// Ajout d'un inconvenient, partie Amir finie\n//\n// 1) Add a new bulleted \"Inconv\u00e9nients\" list item right after the paragraph\n//    that ends with \"...pourraient s'av\u00e9rer plus compliques que d'autres.\"\n//    (same ListParagraph style / numId=6 list / font size as its siblings).\n// 2) Drop the stale <w:lastRenderedPageBreak/> that sits in front of the\n//    \"4. Les diagrammes de classes d\u00e9taill\u00e9s\" Heading2 run.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// --- 1) Insert the new list item -----------------------------------------\nconst anchorText =\n  \"Dans notre cas, les tests unitaires pourraient s\\u2019av\\u00e9rer difficiles, \" +\n  \"il est difficile d\\u2019isoler une classe d\\u00e9pendant d\\u2019une classe Singleton, \" +\n  \"donc certains tests unitaires pourraient s\\u2019av\\u00e9rer plus compliques que d\\u2019autres.\";\n\nconst anchorParagraph = items.find((p) => p.text === anchorText);\nif (!anchorParagraph) {\n  throw new Error(\"Could not locate the 'Dans notre cas...' list paragraph.\");\n}\n\nconst newText =\n  \"S\\u2019il y a une mauvaise conception dans notre impl\\u00e9mentation, \" +\n  \"le patron Singleton peut masquer cela, et causer des probl\\u00e8mes plus tard.\";\n\n// insertParagraph(\"...\", \"After\") clones the anchor paragraph's pPr/rPr\n// (pStyle=ListParagraph, numPr ilvl=0/numId=6, spacing, sz/szCs=24), which is\n// exactly the formatting the new bullet needs.\nanchorParagraph.insertParagraph(newText, \"After\");\n\n// --- 2) Remove the orphan lastRenderedPageBreak ---------------------------\nconst headingText = \"4. Les diagrammes de classes d\\u00e9taill\\u00e9s\";\nconst headingParagraph = items.find((p) => p.text === headingText);\nif (!headingParagraph) {\n  throw new Error(\"Could not locate the '4. Les diagrammes...' heading paragraph.\");\n}\n\n// Re-writing the run text in place drops the stale lastRenderedPageBreak\n// marker while leaving every other run/paragraph property untouched.\nconst headingRange = headingParagraph.getRange();\nheadingRange.insertText(headingText, \"Replace\");\n\nawait context.sync();\n", "ps1": "# Ajout d'un inconvenient, partie Amir finie\n#\n# 1) Add a new bulleted \"Inconvenients\" list item right after the paragraph\n#    that ends with \"...pourraient s'averer plus compliques que d'autres.\"\n#    (same ListParagraph style / numId=6 list / font size as its siblings).\n# 2) Drop the stale LastRenderedPageBreak that sits in front of the\n#    \"4. Les diagrammes de classes detailles\" Heading2 run.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"Dans notre cas, les tests unitaires pourraient s\u2019av\u00e9rer difficiles, il est difficile d\u2019isoler une classe d\u00e9pendant d\u2019une classe Singleton, donc certains tests unitaires pourraient s\u2019av\u00e9rer plus compliques que d\u2019autres.\"\n$headingText = \"4. Les diagrammes de classes d\u00e9taill\u00e9s\"\n\n$anchorParagraph = $null\n$headingParagraph = $null\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.TrimEnd(\"`r\", \"`a\", \"`v\")\n    if ($text -eq $anchorText) {\n        $anchorParagraph = $p\n    }\n    if ($text -eq $headingText) {\n        $headingParagraph = $p\n    }\n}\n\nif ($anchorParagraph -eq $null) {\n    throw \"Could not locate the 'Dans notre cas...' list paragraph.\"\n}\nif ($headingParagraph -eq $null) {\n    throw \"Could not locate the '4. Les diagrammes...' heading paragraph.\"\n}\n\n# --- 1) Remove the orphan LastRenderedPageBreak ----------------------------\n# Re-writing the run text in place drops the stale page-break marker while\n# leaving every other run/paragraph property untouched. Do this BEFORE the\n# structural insert below so the $headingParagraph reference is still live\n# (inserting a new paragraph elsewhere in the story can re-seat/stale out\n# paragraph references captured earlier in the same batch).\n$headingParagraph.Range.Text = $headingText\n\n# --- 2) Insert the new list item -------------------------------------------\n# InsertParagraphAfter() clones the anchor paragraph's formatting (pStyle =\n# ListParagraph, numPr ilvl=0/numId=6, spacing, sz/szCs=24) onto the freshly\n# minted paragraph; we only need to fill in its text afterwards.\n$anchorParagraph.Range.InsertParagraphAfter()\n$newParagraph = $anchorParagraph.Next()\n$newParagraph.Range.Text = \"S\u2019il y a une mauvaise conception dans notre impl\u00e9mentation, le patron Singleton peut masquer cela, et causer des probl\u00e8mes plus tard.\"\n\n\"done\"\n"}
